$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The workbook rows 2, 4 and 5 get their data rotated:
#   New Row2 <- Old Row5 (D, J, K, L, M, P values; I/N/Q unchanged)
#   New Row4 <- Old Row2 (D, J, K, L, M, I, N, P, Q values)
#   New Row5 <- Old Row4 (D, J, K, L, M, I, N, P, Q values)

# --- Row 2 ---
$ws.Range("D2").Value = 44377
$ws.Range("J2").Value = 550
$ws.Range("K2").Value = 2000
$ws.Range("L2").Value = 2800
$ws.Range("M2").Value = 2364
$ws.Range("P2").Value = 394

# --- Row 4 ---
$ws.Range("D4").Value = 44370
$ws.Range("I4").Value = "Segunda"
$ws.Range("J4").Value = 100
$ws.Range("K4").Value = 1000
$ws.Range("L4").Value = 1200
$ws.Range("M4").Value = 1080
$ws.Range("N4").Value = "$/docena de matas"
$ws.Range("P4").Value = 180
$ws.Range("Q4").Value = 6

# --- Row 5 ---
$ws.Range("D5").Value = 44623
$ws.Range("I5").Value = "Primera"
$ws.Range("J5").Value = 300
$ws.Range("K5").Value = 1800
$ws.Range("L5").Value = 2000
$ws.Range("M5").Value = 1900
$ws.Range("N5").Value = "$/paquete"
$ws.Range("P5").Value = 1900
$ws.Range("Q5").Value = 1
